# Add 7 new rows (177-183) of landscaping data for 6/4/2025 (serial 45812),
# continuing the pattern already present in the sheet, and extend the
# Temp_Diff (column F) formula down into the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 177
$lastNewRow  = 183

# 1) Prime the new rows with the same formatting as the last existing row
#    (176) so the date column picks up the existing built-in date style
#    instead of minting a new one.
$ws.Range("A176:T176").Copy()
$ws.Range("A$firstNewRow`:T$lastNewRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2) Row data, taken from the day's repeating 7-row pattern
#    (Flowering/Large, Nonflowering/Medium, Nonflowering/Small,
#     Nonflowering/Medium, Nonflowering/Medium, Nonflowering/Large,
#     Tree/Medium), all for Low=66 / High=87 on 6/4/2025.
$rows = @(
    @{ Row=177; A=45812; B="Flowering";     C="Large";  D=66; E=87; G=0; H=0.1;  I="No"; J=2; K="Bright";  L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 },
    @{ Row=178; A=45812; B="Nonflowering";  C="Medium"; D=66; E=87; G=0; H=0.1;  I="No"; J=3; K="Bright";  L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 },
    @{ Row=179; A=45812; B="Nonflowering";  C="Small";  D=66; E=87; G=0; H=0.2;  I="No"; J=3; K="Bright";  L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 },
    @{ Row=180; A=45812; B="Nonflowering";  C="Medium"; D=66; E=87; G=0; H=0.3;  I="No"; J=3; K="Neutral"; L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 },
    @{ Row=181; A=45812; B="Nonflowering";  C="Medium"; D=66; E=87; G=0; H=0.25; I="No"; J=3; K="Neutral"; L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 },
    @{ Row=182; A=45812; B="Nonflowering";  C="Large";  D=66; E=87; G=0; H=0;    I="No"; J=4; K="Bright";  L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 },
    @{ Row=183; A=45812; B="Tree";          C="Medium"; D=66; E=87; G=0; H=0.9;  I="No"; J=1; K="Neutral"; L=7; M=0.48; N=64; O=30.14; P=18; Q=0.43; R=8.7; S=70; T=38 }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
}

# 3) Extend the Temp_Diff formula (=ABS(Low-High)) down through the new rows.
$ws.Range("F$firstNewRow`:F$lastNewRow").Formula = "=ABS(D$firstNewRow-E$firstNewRow)"

# 4) Update the sheet selection to reflect where the user ended up after
#    adding these rows (the newly-filled Pressure column).
$ws.Range("O177:O183").Select()
